# "All files created with Jan"
# Update the PCB-footprints log entry, add three new rows to the log
# (2022-10-06, 2022-10-07, 2022-10-10) and leave the selection on E33.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Row 30: update wording, and fill in the missing end time ----------
$ws.Range("E30").Value = "Worked on PCB footprints and added 3d models"

$ws.Range("C30").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4122) | Out-Null     # xlPasteFormats -> reuse the h:mm time style
$ws.Range("D30").Value = 0.54166666666666663

# --- Row 31: 2022-10-06, short "x" day ----------------------------------
$ws.Range("B30").Copy() | Out-Null
$ws.Range("B31").PasteSpecial(-4122) | Out-Null     # xlPasteFormats -> reuse the date style
$ws.Range("B31").Value = 44840

$ws.Range("C31").Value = "x"
$ws.Range("D31").Value = "x"
$ws.Range("E31").Value = "x "

# --- Row 32: 2022-10-07 ---------------------------------------------------
$ws.Range("B30").Copy() | Out-Null
$ws.Range("B32").PasteSpecial(-4122) | Out-Null
$ws.Range("B32").Value = 44841

$ws.Range("C30").Copy() | Out-Null
$ws.Range("C32").PasteSpecial(-4122) | Out-Null
$ws.Range("C32").Value = 0.4375

$ws.Range("D30").Copy() | Out-Null
$ws.Range("D32").PasteSpecial(-4122) | Out-Null
$ws.Range("D32").Value = 0.58333333333333337

$ws.Range("E32").Value = "Worked on pcb layout, improved pcb workflow skills (to prevent wasting time)"

# --- Row 33: 2022-10-10 (no end time filled in) ---------------------------
$ws.Range("B30").Copy() | Out-Null
$ws.Range("B33").PasteSpecial(-4122) | Out-Null
$ws.Range("B33").Value = 44844

$ws.Range("C30").Copy() | Out-Null
$ws.Range("C33").PasteSpecial(-4122) | Out-Null
$ws.Range("C33").Value = 0.3888888888888889

$ws.Range("E33").Value = "Worked on software architecture, with the help of Jan."

$excel.CutCopyMode = 0

# --- Leave the view/selection where the author left it --------------------
$ws.Range("E33").Select() | Out-Null

Write-Output "done"
